$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Test case 2.1: "numeric" -> "integer"
$ws.Range("B11").Value = "Incorrect guess: integer. Guesses>0"

# Test case 2.2: "non-numeric" -> "non-integer"
$ws.Range("B14").Value = "incorrect guess non-integer (including blank)"

# Reflect the new selection left by the edit
$ws.Activate()
$ws.Range("B6").Select()
